$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) updates ---
# "Amount" -> "Payment Amount"
$ws.Range("AE1").Value = "Payment Amount"

# New expense-tracking columns appended after Country (AG1)
$ws.Range("AH1").Value = "Expense Type"
$ws.Range("AI1").Value = "Expense Amount"
$ws.Range("AJ1").Value = "Expense Number"

# --- Data row (row 2) updates ---
# S Invoice Ref timestamp moves forward to the new expense-creation run
$ws.Range("C2").Value = "Quick 2019/11/26 18:58:19"

# Invoice Number replaced by the newly generated Expense Number for this row
$ws.Range("AB2").Value = "EXP-526-261119-8"

# Fill in the new expense columns for row 2
$ws.Range("AH2").Value = "Shipping Expense"
$ws.Range("AI2").Value = 100

# --- View state: scroll / selection moved while the new columns were added ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 24
$win.ScrollRow = 1
$ws.Range("AE8").Select()
